# delete dissolved municipalities from options
#
# The "area" sheet lists geographic options (municipalities / stadsdelen /
# wijken) used elsewhere in the workbook. Two municipalities that have since
# been dissolved/merged ("Weesp" and "Beemster") are removed from that list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("area")

# Locate and delete the "Weesp" row (column A holds the geography name).
$weespRow = $ws.Cells.Find("Weesp").Row
$ws.Rows.Item($weespRow).Delete()

# Locate and delete the "Beemster" row (shifts up after the first delete).
$beemsterRow = $ws.Cells.Find("Beemster").Row
$ws.Rows.Item($beemsterRow).Delete()

# Make "area" the active sheet, with the row that now starts with "Zuidoost"
# (originally row 28, now row 27) selected - matching the saved view state.
$ws.Activate()
$ws.Range("A27:B27").Select()
